$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create "Norway" and "Poland" worksheets by cloning an existing, already
#    laid-out "accessories" sheet (Greece) so the column widths / row
#    heights / merged cells / page setup all match without extra fix-up.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("Greece")

$template.Copy($null, $wb.Worksheets.Item("Hungary"))
$norway = $wb.Worksheets.Item($wb.Worksheets.Count)
$norway.Name = "Norway"

$template.Copy($null, $wb.Worksheets.Item("Norway"))
$poland = $wb.Worksheets.Item($wb.Worksheets.Count)
$poland.Name = "Poland"

# ---------------------------------------------------------------------------
# 2) Fill in the market-specific values. B4 (article code) before B2 (market
#    name) on each sheet, Norway fully before Poland -- this reproduces the
#    shared-string insertion order seen in the target file (idx 47-50).
# ---------------------------------------------------------------------------
$norway.Range("B4").Value = "NGC-2931/T3086/T3085"
$norway.Range("B2").Value = "Norway Market"

$poland.Range("B4").Value = "NGC-2920/T3037/T3120"
$poland.Range("B2").Value = "Poland Market"

# ---------------------------------------------------------------------------
# 3) Insert two new accessory rows ("MX-BBX" / "MX-DPBX") above the old row 9
#    on every sheet that needs them: Greece, Croatia, Portugal (reversed
#    order), Norway, Poland.
# ---------------------------------------------------------------------------
function Add-AccessoryRows($ws, $first, $second) {
    $ws.Rows("9:10").Insert()
    $ws.Range("A8").Copy()
    $ws.Range("A9:A10").PasteSpecial(-4122)
    $ws.Range("A9").Value = $first
    $ws.Range("A10").Value = $second
    $ws.Range("A10").Select()
}

Add-AccessoryRows $wb.Worksheets.Item("Greece")   "MX-BBX"  "MX-DPBX"
Add-AccessoryRows $wb.Worksheets.Item("Croatia")  "MX-BBX"  "MX-DPBX"
Add-AccessoryRows $norway                          "MX-BBX"  "MX-DPBX"
Add-AccessoryRows $poland                           "MX-BBX"  "MX-DPBX"
Add-AccessoryRows $wb.Worksheets.Item("Portugal") "MX-DPBX" "MX-BBX"

$poland.Range("A9").Select()
$poland.Activate()
